$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared/base formulas so the output file extension changes
# from ".s" to ".asm". Excel will propagate the edit across each shared
# formula group (B1, B2:B65, B66:B77) and recompute the cached <v> values.
$ws.Range("B1").Formula = '=CONCATENATE("riscv32-unknown-elf-objdump -d test/", A1, ".riscv"," >  test/", A1, ".asm")'
$ws.Range("B2:B65").Formula = '=CONCATENATE("riscv32-unknown-elf-objdump -d test/", A2, ".riscv"," >  test/", A2, ".asm")'
$ws.Range("B66:B77").Formula = '=CONCATENATE("riscv32-unknown-elf-objdump -d test/", A66, ".riscv"," >  test/", A66, ".asm")'

# Update the sheet selection to cover the whole column B (B1:B1048576)
# instead of just B1:B77, matching the saved view state.
$ws.Range("B1:B1048576").Select()
